# Estructuracion de reporte serenityReport
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the license-plate value in A2 ("WFL326" -> "AMZ45D")
$ws.Range("A2").Value = "AMZ45D"

# Apply the new Arial/black font to A2
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Color = 0

# Move the active selection from C6 to A2
$ws.Range("A2").Select()
